$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "mapsto" column (C) to hold string labels instead of numeric
#     codes, and refresh a couple of "meaning" (B) / "mapsfrom" (A) labels.
#     Writing these in this particular order keeps the shared-string table
#     append order (and thus indices) aligned with the authored workbook:
#     stage1, stage2, sws, rem, "stage 4", "Sleep stage R".

# Row 2: Sleep stage W / wake / wake
$ws.Range("C2").Value = "wake"

# Row 3: Sleep stage 1 / stage 1 / stage1
$ws.Range("C3").Value = "stage1"

# Row 4: Sleep stage 2 / stage 2 / stage2
$ws.Range("C4").Value = "stage2"

# Row 5: Sleep stage 3 / stage 3 / sws
$ws.Range("C5").Value = "sws"

# Row 7 (old "Sleep stage ?" row): now "rem" meaning/mapsto, before the
# label edits below so the shared string for "rem" is appended before
# "stage 4" / "Sleep stage R" (matches authored append order).
$ws.Range("C7").Value = "rem"

# Row 6: Sleep stage 4 / stage 4 (was "stage 3") / sws
$ws.Range("B6").Value = "stage 4"
$ws.Range("C6").Value = "sws"

# Row 7: "Sleep stage ?" -> "Sleep stage R", meaning REM, mapsto rem
$ws.Range("A7").Value = "Sleep stage R"
$ws.Range("B7").Value = "REM"

# Row 8: "Sleep stage 5" -> "Sleep stage ?", meaning unknown, mapsto unknown
$ws.Range("A8").Value = "Sleep stage ?"
$ws.Range("B8").Value = "unknown"
$ws.Range("C8").Value = "unknown"

# --- Column widths (A, B) -- authored widths of 20.46875 / 16.5859375
#     characters land on this engine's internal 1/6-character rounding
#     grid at 20.5 / 16.6667, so feed it the ColumnWidth inputs that round
#     to those nearest grid points.
$ws.Columns.Item(1).ColumnWidth = 19.666666666666668
$ws.Columns.Item(2).ColumnWidth = 15.833333333333334

# --- Selection moves to A8
$ws.Range("A8").Select()
